$wb = $excel.ActiveWorkbook

# The three new "Market" sheets are modeled directly on the existing
# "Denmark" sheet (same layout / merged cells / styles), so build each
# one by copying Denmark (and then copying forward) and only touching
# the handful of cells + cosmetic bits that differ per-country.
$denmark = $wb.Worksheets.Item("Denmark")

$denmark.Copy($null, $denmark)
$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"

$russia.Copy($null, $russia)
$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"

$finland.Copy($null, $finland)
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"

# --- Russia --- (ticket code entered before the market name)
$russia.Range("B4").Value = "NGC-2929/T2925"
$russia.Range("B2").Value = "Russia Market"
$russia.Rows.Item(3).RowHeight = 28.8
$russia.Rows.Item(4).RowHeight = 28.8
$russia.Rows.Item(5).RowHeight = 28.8
$russia.Columns.Item(2).ColumnWidth = 15.21875
$russia.Columns.Item(3).ColumnWidth = 8.43
$russia.Columns.Item(4).ColumnWidth = 8.44140625
$russia.Range("A1:D10").Select()

# --- Finland ---
$finland.Range("B4").Value = "NGC-3130/T2957"
$finland.Range("B2").Value = "Finland Market"
$finland.Rows.Item(3).RowHeight = 28.8
$finland.Rows.Item(4).RowHeight = 28.8
$finland.Rows.Item(5).RowHeight = 28.8
$finland.Columns.Item(2).ColumnWidth = 15.21875
$finland.Columns.Item(3).ColumnWidth = 8.43
$finland.Columns.Item(4).ColumnWidth = 8.44140625
$finland.Range("A1:D10").Select()

# --- Hungary ---
$hungary.Range("B4").Value = "NGC-3104/T2979"
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Rows.Item(3).RowHeight = 28.8
$hungary.Rows.Item(4).RowHeight = 28.8
$hungary.Rows.Item(5).RowHeight = 28.8
$hungary.Columns.Item(2).ColumnWidth = 15.21875
$hungary.Columns.Item(3).ColumnWidth = 8.43
$hungary.Columns.Item(4).ColumnWidth = 8.44140625

# Hungary ends up as the new active/visible tab, with H17 selected.
$hungary.Activate()
$hungary.Range("H17").Select()
